# Weekly data refresh: a new week's price record is inserted at row 587
# (right after the header + prior records), pushing every existing
# record down by one row. The sheet's used range grows from A1:T690 to
# A1:T691.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 587; this shifts rows
# 587..690 down to 588..691 (so the former row 690 becomes row 691).
$ws.Rows.Item(587).Insert()

# Populate the newly inserted row 587 with the new week's record.
$ws.Cells.Item(587, 1).Value2  = 4
$ws.Cells.Item(587, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(587, 3).Value2  = "Los Lagos"
$ws.Cells.Item(587, 4).Value2  = 45218
$ws.Cells.Item(587, 5).Value2  = 10
$ws.Cells.Item(587, 6).Value2  = "Fruta"
$ws.Cells.Item(587, 7).Value2  = 100102
$ws.Cells.Item(587, 8).Value2  = "Cítricos"
$ws.Cells.Item(587, 9).Value2  = 100102006
$ws.Cells.Item(587, 10).Value2 = "Pomelo"
$ws.Cells.Item(587, 11).Value2 = "Start Ruby"
$ws.Cells.Item(587, 12).Value2 = "Primera"
$ws.Cells.Item(587, 13).Value2 = 100
$ws.Cells.Item(587, 14).Value2 = 15000
$ws.Cells.Item(587, 15).Value2 = 15000
$ws.Cells.Item(587, 16).Value2 = 15000
$ws.Cells.Item(587, 17).Value2 = "$/caja 14 kilos empedrada"
$ws.Cells.Item(587, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(587, 19).Value2 = 1071
$ws.Cells.Item(587, 20).Value2 = 14
